$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells so numeric-looking strings
# (e.g. "1.400", "10.000") keep their exact literal text like the source data.
foreach ($addr in @("D2","D3","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D22","D25","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.009.48"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "1.651.85"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "216.58"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").Value = "0.5202"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "0.2620"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").Value = "0.06259"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").Value = "20.57"
$ws.Range("E10").Value = "  -3.89%  "

$ws.Range("D11").Value = "0.07713"
$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.688.80"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.470"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "1.878.43"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "0.5421"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").Value = "0.0₅8084"

$ws.Range("D17").Value = "64.76"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").Value = "26.022.56"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("E20").Value = "  -2.60%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "10.000"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("E23").Value = "  -3.58%  "

$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").Value = "138.07"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("D27").Value = "7.239"
$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("D28").Value = "16.05"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("D29").Value = "1.400"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "0.05933"
$ws.Range("E30").Value = "  -2.08%  "

$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("D32").Value = "3.512"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").Value = "3.244"
$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("E34").Value = "  -5.99%  "

$ws.Range("D35").Value = "0.9486"
$ws.Range("E35").Value = "  -3.91%  "

$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("D37").Value = "2.754"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("D38").Value = "0.5675"
$ws.Range("E38").Value = "  -4.75%  "

$ws.Range("D39").Value = "0.01594"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").Value = "5.882"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("D41").Value = "0.8457"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "100.56"
$ws.Range("E43").Value = "  +0.63%  "

$ws.Range("D44").Value = "999.00"
$ws.Range("E44").Value = "  -6.21%  "

$ws.Range("D45").Value = "1.792.66"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "56.59"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "7.969"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").Value = "0.4299"
$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("E51").Value = "  +0.02%  "
